$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.191.99"
$ws.Range("E2").Value = "  -1.71%  "

$ws.Range("D3").Value = "'2.273.34"
$ws.Range("E3").Value = "  -1.76%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'111.72"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").Value = "'265.10"
$ws.Range("E6").Value = "  -2.62%  "

$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  -3.59%  "

$ws.Range("D10").Value = "'47.74"
$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").Value = "'0.0928"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").Value = "'8.91"
$ws.Range("E12").Value = "  +1.49%  "

$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").Value = "'15.38"
$ws.Range("E14").Value = "  -2.85%  "

$ws.Range("D15").Value = "'2.613.37"

$ws.Range("D16").Value = "'0.853"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("D17").Value = "'2.280.31"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").Value = "'43.095.99"
$ws.Range("E18").Value = "  -1.82%  "

$ws.Range("D19").Value = "'0.0000108"
$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("D20").Value = "'6.86"
$ws.Range("E20").Value = "  +4.74%  "

$ws.Range("D21").Value = "'71.16"
$ws.Range("E21").Value = "  -2.05%  "

$ws.Range("D22").Value = "'2.46"
$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("D23").Value = "'231.40"

$ws.Range("D24").Value = "'9.70"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("D25").Value = "'2.87"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'11.35"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("D28").Value = "'3.92"
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").Value = "'40.35"
$ws.Range("E29").Value = "  -6.43%  "

$ws.Range("D30").Value = "'3.35"
$ws.Range("E30").Value = "  -3.25%  "

$ws.Range("D31").Value = "'2.25"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").Value = "'171.40"
$ws.Range("E32").Value = "  -3.93%  "

$ws.Range("D33").Value = "'21.34"
$ws.Range("E33").Value = "  -2.50%  "

$ws.Range("D34").Value = "'0.0904"
$ws.Range("E34").Value = "  -4.38%  "

$ws.Range("D35").Value = "'5.79"
$ws.Range("E35").Value = "  +2.84%  "

$ws.Range("D36").Value = "'0.127"
$ws.Range("E36").Value = "  -0.68%  "

$ws.Range("D37").Value = "'4.67"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").Value = "'0.0351"
$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("D39").Value = "'3.85"
$ws.Range("E39").Value = "  -2.13%  "

$ws.Range("D40").Value = "'0.105"
$ws.Range("E40").Value = "  -7.08%  "

$ws.Range("D41").Value = "'2.57"
$ws.Range("E41").Value = "  +7.17%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'14.13"
$ws.Range("E42").Value = "  +14.69%  "

$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'75.12"
$ws.Range("E43").Value = "  +6.80%  "

$ws.Range("D44").Value = "'0.237"
$ws.Range("E44").Value = "  -2.29%  "

$ws.Range("D45").Value = "'6.09"
$ws.Range("E45").Value = "  +9.65%  "

$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "'1.37"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").Value = "'8.69"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").Value = "'0.0991"
$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.24"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'100.67"
$ws.Range("E51").Value = "  +0.65%  "
